$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix formatting on fastq purpose column: "fullRNASEQ" -> "fullRNASeq"
for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # Column E = purpose
    if ($cell.Value2 -ceq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
